$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New density value (column D) applies to every data row (2-51)
$newDensity = 0.3379310344827586

# New "steps" values (column E) for rows 2-51, in order
$stepsValues = @(7,11,7,11,9,6,6,10,8,7,7,8,9,7,8,12,11,8,6,10,8,8,7,7,7,8,13,7,8,8,7,13,8,6,9,7,5,8,9,10,7,8,8,7,12,7,6,14,7,9)

for ($i = 0; $i -lt $stepsValues.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 4).Value = $newDensity
    $ws.Cells.Item($row, 5).Value = $stepsValues[$i]
}
